$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,("ECs","Fgf2","Fgfr1","ECs",1,0.3333333333333333,0.2151743333333333,0.645523,0.01945888736810517,0.01945888736810517,3,1,3.76917,11.30751,0.02686645020528053,0.02686645020528053,0.81102864197,7.29925777773,0.0005227912285253598,0.0005227912285253598)
    ,("ECs","Fgf2","Fgfr1","FAPs",1,0.3333333333333333,0.2151743333333333,0.645523,0.01945888736810517,0.01945888736810517,3,1,82.48060333333333,247.44181,0.5879175050094569,0.5879175050094571,17.74770883518111,159.72937951663,0.01144022051171643,0.01144022051171643)
    ,("ECs","Fgf2","Fgfr1","M2",1,0.3333333333333333,0.2151743333333333,0.645523,0.01945888736810517,0.01945888736810517,3,1,0.5716056666666667,1.714817,0.004074375838860061,0.004074375838860062,0.1229948682545556,1.106953814291,0.00007928282054370695,0.00007928282054370696)
    ,("ECs","Fgf2","Fgfr1","sCs",1,0.3333333333333333,0.2151743333333333,0.645523,0.01945888736810517,0.01945888736810517,3,1,53.47143866666666,160.414316,0.3811416689464024,0.3811416689464024,11.50568116747422,103.551130507268,0.007416592807319672,0.007416592807319673)
    ,("FAPs","Fgf2","Fgfr1","ECs",3,1,8.077181333333334,24.231544,0.730444748601188,0.730444748601188,3,1,3.76917,11.30751,0.02686645020528053,0.02686645020528053,30.44426956616001,273.99842609544,0.01962445746600247,0.01962445746600247)
    ,("FAPs","Fgf2","Fgfr1","FAPs",3,1,8.077181333333334,24.231544,0.730444748601188,0.730444748601188,3,1,82.48060333333333,247.44181,0.5879175050094569,0.5879175050094571,666.2107896060711,5995.89710645464,0.4294412541448705,0.4294412541448706)
    ,("FAPs","Fgf2","Fgfr1","M2",3,1,8.077181333333334,24.231544,0.730444748601188,0.730444748601188,3,1,0.5716056666666667,1.714817,0.004074375838860061,0.004074375838860062,4.616962620827556,41.552663587448,0.002976106435322892,0.002976106435322892)
    ,("FAPs","Fgf2","Fgfr1","sCs",3,1,8.077181333333334,24.231544,0.730444748601188,0.730444748601188,3,1,53.47143866666666,160.414316,0.3811416689464024,0.3811416689464024,431.8985062648782,3887.086556383903,0.2784029305549921,0.2784029305549922)
    ,("M2","Fgf2","Fgfr1","ECs",3,1,0.519176,1.557528,0.04695070806875992,0.04695070806875992,3,1,3.76917,11.30751,0.02686645020528053,0.02686645020528053,1.95686260392,17.61176343528,0.001261398860432001,0.001261398860432001)
    ,("M2","Fgf2","Fgfr1","FAPs",3,1,0.519176,1.557528,0.04695070806875992,0.04695070806875992,3,1,82.48060333333333,247.44181,0.5879175050094569,0.5879175050094571,42.82194971618667,385.39754744568,0.02760314314621271,0.02760314314621272)
    ,("M2","Fgf2","Fgfr1","M2",3,1,0.519176,1.557528,0.04695070806875992,0.04695070806875992,3,1,0.5716056666666667,1.714817,0.004074375838860061,0.004074375838860062,0.2967639435973333,2.670875492376,0.0001912948305727275,0.0001912948305727276)
    ,("M2","Fgf2","Fgfr1","sCs",3,1,0.519176,1.557528,0.04695070806875992,0.04695070806875992,3,1,53.47143866666666,160.414316,0.3811416689464024,0.3811416689464024,27.76108764120533,249.849788770848,0.01789487123154248,0.01789487123154248)
    ,("sCs","Fgf2","Fgfr1","ECs",3,1,2.246363333333333,6.73909,0.2031456559619469,0.2031456559619469,3,1,3.76917,11.30751,0.02686645020528053,0.02686645020528053,8.466925285100002,76.20232756590001,0.005457802650320697,0.005457802650320697)
    ,("sCs","Fgf2","Fgfr1","FAPs",3,1,2.246363333333333,6.73909,0.2031456559619469,0.2031456559619469,3,1,82.48060333333333,247.44181,0.5879175050094569,0.5879175050094571,185.2814030392111,1667.5326273529,0.1194328872066574,0.1194328872066574)
    ,("sCs","Fgf2","Fgfr1","M2",3,1,2.246363333333333,6.73909,0.2031456559619469,0.2031456559619469,3,1,0.5716056666666667,1.714817,0.004074375838860061,0.004074375838860062,1.284034010725556,11.55630609653,0.000827691752420735,0.0008276917524207352)
    ,("sCs","Fgf2","Fgfr1","sCs",3,1,2.246363333333333,6.73909,0.2031456559619469,0.2031456559619469,3,1,53.47143866666666,160.414316,0.3811416689464024,0.3811416689464024,120.1162792013822,1081.04651281244,0.07742727435254813,0.07742727435254815)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
}
